# #5: property boat&car done
# Add a header row (with a new "capacity" column) and extra metadata
# columns to the 汽車 (car) worksheet, matching the other property sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車

# --- Row 1: proper column headers (previously row 1 just duplicated the
#     data row with no header text) -----------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: correct the car data and add the trailing metadata columns
$ws.Range("B2").Value = "toyotarav4rod"
$ws.Range("E2").Value = "100年03月01曰"
$ws.Range("F2").Value = "買賣"

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-30"
$ws.Range("K2").Value = "林淑芬"
$ws.Range("L2").Value = 1337
$ws.Range("M2").Value = "tmp63cf1"
$ws.Range("N2").Value = 29
